$wb = $excel.ActiveWorkbook

# Row 2 on sheet ALC (context G value 5489)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 98.666664
$ws.Range("I2").Value = 98.666664
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 98.666664
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 14.333336
$ws.Range("N2").ClearContents()

# Row 21 on sheet ALC (context G value 2149)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 27732.23
$ws.Range("J21").Value = 28901.9
$ws.Range("L21").Value = 28901.9
$ws.Range("N21").Value = -29837.9

# Row 23 on sheet ALC (context G value 2149)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 27732.23
$ws.Range("J23").Value = 28901.9
$ws.Range("L23").Value = 28901.9
$ws.Range("N23").Value = -29369.9

# Row 29 on sheet ALC (context G value 4575)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 93.75
$ws.Range("I29").Value = 93.75
$ws.Range("K29").Value = 281.25
$ws.Range("M29").Value = -0.25

# Row 38 on sheet ALC (context G value 4599)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1493.3478
$ws.Range("I38").Value = 77.09090999999999
$ws.Range("J38").Value = 2791.5833
$ws.Range("K38").Value = 231.27273
$ws.Range("L38").Value = 8374.749899999999
$ws.Range("M38").Value = 140.72727
$ws.Range("N38").Value = -9118.749899999999

# Row 55 on sheet ALC (context G value 5517)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 173.33333
$ws.Range("I55").Value = 110
$ws.Range("K55").Value = 110
$ws.Range("M55").Value = 104

# Row 64 on sheet ALC (context G value 5506)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3197.6572
$ws.Range("I64").Value = 2938.6875
$ws.Range("J64").Value = 3415.7368
$ws.Range("K64").Value = 2938.6875
$ws.Range("L64").Value = 3415.7368
$ws.Range("M64").Value = -2690.6875
$ws.Range("N64").Value = -3911.7368

# Row 67 on sheet ALC (context G value 5506)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3197.6572
$ws.Range("I67").Value = 2938.6875
$ws.Range("J67").Value = 3415.7368
$ws.Range("K67").Value = 2938.6875
$ws.Range("L67").Value = 3415.7368
$ws.Range("M67").Value = -2080.6875
$ws.Range("N67").Value = -5131.736800000001

# Row 92 on sheet ALC (context G value 19901)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 596.125
$ws.Range("I92").Value = 217.25
$ws.Range("K92").Value = 217.25
$ws.Range("M92").Value = 1030.75

# Row 96 on sheet ALC (context G value 19894)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 466.58334
$ws.Range("I96").Value = 448.33334
$ws.Range("K96").Value = 1345.00002
$ws.Range("M96").Value = 27.99998000000005

# Row 100 on sheet ALC (context G value 19906)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2867.2222
$ws.Range("I100").Value = 1250
$ws.Range("J100").Value = 3329.2856
$ws.Range("K100").Value = 1250
$ws.Range("L100").Value = 3329.2856
$ws.Range("M100").Value = -709
$ws.Range("N100").Value = -4411.2856

# Row 123 on sheet ALC (context G value 34090)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 69791
$ws.Range("J123").Value = 69791
$ws.Range("L123").Value = 69791
$ws.Range("N123").Value = -79591

# Row 141 on sheet ALC (context G value 44161)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 25745.834
$ws.Range("I141").Value = 12618.75
$ws.Range("J141").Value = 52000
$ws.Range("K141").Value = 37856.25
$ws.Range("L141").Value = 156000
$ws.Range("M141").Value = -32676.25
$ws.Range("N141").Value = -166360

# Row 4 on sheet ARM (context G value 5071)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 640.61536
$ws.Range("I4").Value = 640.61536
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 640.61536
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -524.61536
$ws.Range("N4").ClearContents()

# Row 101 on sheet ARM (context G value 18518)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 53600.8
$ws.Range("J101").Value = 53600.8
$ws.Range("L101").Value = 53600.8
$ws.Range("N101").Value = -60090.8

# Row 110 on sheet ARM (context G value 27708)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2167.95
$ws.Range("I110").Value = 2236.6667
$ws.Range("J110").Value = 1549.5
$ws.Range("K110").Value = 2236.6667
$ws.Range("L110").Value = 1549.5
$ws.Range("M110").Value = -191.6667000000002
$ws.Range("N110").Value = -5639.5

# Row 135 on sheet ARM (context G value 42016)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").ClearContents()

# Row 94 on sheet BSM (context G value 19939)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3441.6667
$ws.Range("I94").Value = 1275
$ws.Range("J94").Value = 4525
$ws.Range("K94").Value = 1275
$ws.Range("L94").Value = 4525
$ws.Range("M94").Value = -824
$ws.Range("N94").Value = -5427

# Row 99 on sheet BSM (context G value 19943)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1781.3334
$ws.Range("I99").Value = 1699.1666
$ws.Range("J99").Value = 2110
$ws.Range("K99").Value = 1699.1666
$ws.Range("L99").Value = 2110
$ws.Range("M99").Value = -201.1666
$ws.Range("N99").Value = -5106

# Row 6 on sheet CUL (context G value 4639)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 95.57143000000001
$ws.Range("I6").Value = 78.166664
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 234.499992
$ws.Range("L6").Value = 600
$ws.Range("M6").Value = -121.499992
$ws.Range("N6").Value = -826

# Row 113 on sheet CUL (context G value 27843)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 417299.62
$ws.Range("I113").Value = 556169.9
$ws.Range("J113").Value = 688.8333
$ws.Range("K113").Value = 1668509.7
$ws.Range("L113").Value = 2066.4999
$ws.Range("M113").Value = -1666339.7
$ws.Range("N113").Value = -6406.4999

# Row 117 on sheet CUL (context G value 27870)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1497.0588
$ws.Range("J117").Value = 1722.3846
$ws.Range("L117").Value = 5167.1538
$ws.Range("N117").Value = -12051.1538

# Row 121 on sheet CUL (context G value 27878)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 917.7568
$ws.Range("I121").Value = 592.375
$ws.Range("J121").Value = 1007.5172
$ws.Range("K121").Value = 1777.125
$ws.Range("L121").Value = 3022.5516
$ws.Range("M121").Value = -467.125
$ws.Range("N121").Value = -5642.5516

# Row 134 on sheet CUL (context G value 44074)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 4608.6484
$ws.Range("I134").Value = 1424.6154
$ws.Range("J134").Value = 6333.3335
$ws.Range("K134").Value = 4273.8462
$ws.Range("L134").Value = 19000.0005
$ws.Range("M134").Value = 796.1538
$ws.Range("N134").Value = -29140.0005

# Row 140 on sheet CUL (context G value 44097)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2043.0968
$ws.Range("I140").Value = 686.4706
$ws.Range("J140").Value = 3690.4285
$ws.Range("K140").Value = 2059.4118
$ws.Range("L140").Value = 11071.2855
$ws.Range("M140").Value = 3120.5882
$ws.Range("N140").Value = -21431.2855

# Row 70 on sheet GSM (context G value 14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8110.9
$ws.Range("J70").Value = 8434.833000000001
$ws.Range("L70").Value = 8434.833000000001
$ws.Range("N70").Value = -8974.833000000001

# Row 73 on sheet GSM (context G value 14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 8110.9
$ws.Range("J73").Value = 8434.833000000001
$ws.Range("L73").Value = 8434.833000000001
$ws.Range("N73").Value = -10306.833

# Row 47 on sheet LTW (context G value 3138)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 25500
$ws.Range("J47").Value = 25500
$ws.Range("L47").Value = 25500
$ws.Range("N47").Value = -26480

# Row 52 on sheet LTW (context G value 3138)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H52").Value = 25500
$ws.Range("J52").Value = 25500
$ws.Range("L52").Value = 25500
$ws.Range("N52").Value = -25966

# Row 93 on sheet LTW (context G value 19993)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1371.4286
$ws.Range("J93").Value = 1380
$ws.Range("L93").Value = 1380
$ws.Range("N93").Value = -3876

# Row 122 on sheet LTW (context G value 36247)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6472.0713
$ws.Range("I122").Value = 5686.2856
$ws.Range("J122").Value = 7257.857
$ws.Range("K122").Value = 17058.8568
$ws.Range("L122").Value = 21773.571
$ws.Range("M122").Value = -14608.8568
$ws.Range("N122").Value = -26673.571

# Row 100 on sheet WVR (context G value 19981)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 11202.35
$ws.Range("I100").Value = 19076.455
$ws.Range("J100").Value = 1578.4445
$ws.Range("K100").Value = 38152.91
$ws.Range("L100").Value = 3156.889
$ws.Range("M100").Value = -37611.91
$ws.Range("N100").Value = -4238.889
